$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The cryptos table stores Price (column D) and Volume(1h) (column E) as plain
# text so values like "22.90" or "1.000" keep their exact digits/trailing zeros.
# Excel auto-converts numeric-looking text typed/assigned into a cell into a real
# number (dropping trailing zeros, switching to scientific notation, etc.), so the
# Price cells that look like plain numbers are pre-formatted as Text before their
# new values are written; this keeps them stored as text exactly as scraped.
$ws.Range('D4:D11,D13:D17,D19:D36,D38:D43,D45:D50').NumberFormat = '@'

$ws.Range('D2').Value = '29.062.10'
$ws.Range('E2').Value = '  +0.08%  '
$ws.Range('D3').Value = '1.834.94'
$ws.Range('E3').Value = '  +0.34%  '
$ws.Range('D4').Value = '0.9999'
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').Value = '244.45'
$ws.Range('E5').Value = '  +1.51%  '
$ws.Range('D6').Value = '0.6343'
$ws.Range('E6').Value = '  +1.89%  '
$ws.Range('D7').Value = '1.001'
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').Value = '0.07565'
$ws.Range('E8').Value = '  +2.00%  '
$ws.Range('D9').Value = '0.2947'
$ws.Range('E9').Value = '  +0.96%  '
$ws.Range('D10').Value = '22.90'
$ws.Range('E10').Value = '  +0.98%  '
$ws.Range('D11').Value = '0.07745'
$ws.Range('E11').Value = '  +1.86%  '
$ws.Range('D12').Value = '1.853.89'
$ws.Range('E12').Value = '  +1.35%  '
$ws.Range('D13').Value = '5.003'
$ws.Range('E13').Value = '  +0.91%  '
$ws.Range('D14').Value = '0.6713'
$ws.Range('D15').Value = '83.25'
$ws.Range('E15').Value = '  +1.49%  '
$ws.Range('D16').Value = '0.000009613'
$ws.Range('E16').Value = '  +5.11%  '
$ws.Range('D17').Value = '6.117'
$ws.Range('E17').Value = '  +1.80%  '
$ws.Range('D18').Value = '29.099.06'
$ws.Range('E18').Value = '  +0.19%  '
$ws.Range('D19').Value = '12.58'
$ws.Range('E19').Value = '  +1.88%  '
$ws.Range('D20').Value = '226.70'
$ws.Range('E20').Value = '  +0.70%  '
$ws.Range('D21').Value = '1.000'
$ws.Range('E21').Value = '  -0.01%  '
$ws.Range('D22').Value = '7.222'
$ws.Range('E22').Value = '  +0.67%  '
$ws.Range('D23').Value = '1.001'
$ws.Range('E23').Value = '  +0.08%  '
$ws.Range('D24').Value = '160.65'
$ws.Range('E24').Value = '  +0.83%  '
$ws.Range('D25').Value = '0.1403'
$ws.Range('E25').Value = '  +3.49%  '
$ws.Range('D26').Value = '8.546'
$ws.Range('E26').Value = '  +1.61%  '
$ws.Range('D27').Value = '17.97'
$ws.Range('E27').Value = '  +1.00%  '
$ws.Range('D28').Value = '1.499'
$ws.Range('E28').Value = '  +0.18%  '
$ws.Range('D29').Value = '4.125'
$ws.Range('E29').Value = '  +1.75%  '
$ws.Range('D30').Value = '4.080'
$ws.Range('E30').Value = '  +1.37%  '
$ws.Range('D31').Value = '1.204'
$ws.Range('E31').Value = '  +0.18%  '
$ws.Range('D32').Value = '0.05426'
$ws.Range('E32').Value = '  +3.58%  '
$ws.Range('D33').Value = '1.863'
$ws.Range('E33').Value = '  +1.49%  '
$ws.Range('D34').Value = '0.7478'
$ws.Range('E34').Value = '  +1.80%  '
$ws.Range('D35').Value = '1.142'
$ws.Range('E35').Value = '  -0.69%  '
$ws.Range('D36').Value = '2.660'
$ws.Range('E36').Value = '  +0.49%  '
$ws.Range('D37').Value = '1.237.26'
$ws.Range('E37').Value = '  -3.37%  '
$ws.Range('D38').Value = '2.761'
$ws.Range('E38').Value = '  +0.42%  '
$ws.Range('D39').Value = '0.01791'
$ws.Range('E39').Value = '  +0.56%  '
$ws.Range('D40').Value = '6.618'
$ws.Range('E40').Value = '  +4.71%  '
$ws.Range('D41').Value = '0.9050'
$ws.Range('E41').Value = '  +1.18%  '
$ws.Range('D42').Value = '1.001'
$ws.Range('E42').Value = '  +0.06%  '
$ws.Range('D43').Value = '102.28'
$ws.Range('D44').Value = '1.996.33'
$ws.Range('E44').Value = '  +0.93%  '
$ws.Range('D45').Value = '0.00000000125'
$ws.Range('E45').Value = '  +4.38%  '
$ws.Range('D46').Value = '65.11'
$ws.Range('E46').Value = '  +2.32%  '
$ws.Range('D47').Value = '0.5111'
$ws.Range('E47').Value = '  -0.06%  '
$ws.Range('D48').Value = '0.4091'
$ws.Range('E48').Value = '  +3.22%  '
$ws.Range('D49').Value = '9.087'
$ws.Range('E49').Value = '  +3.26%  '
$ws.Range('D50').Value = '6.774'
$ws.Range('E50').Value = '  +1.74%  '
$ws.Range('E51').Value = '  +0.54%  '
